# Update "Data Analyst project tracker.xlsx":
# Fill in rows 8-10 (PowerBI learning activities) that were previously only
# partially populated (SL#, Topic, Activity were present; Date, Effort,
# Status, Remarks were still blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108
$xlTop = -4160

# --- Row 8: "Connect PBI to DB and pull data" ---------------------------
$ws.Range("D8").NumberFormat = "d-mmm-yy"
$ws.Range("D8").Value = 44305
$ws.Range("E8").Value = "2hr"
$ws.Range("E8").HorizontalAlignment = $xlCenter
$ws.Range("E8").VerticalAlignment = $xlCenter
$ws.Range("F8").Value = "completed"
$ws.Range("F8").HorizontalAlignment = $xlCenter

# --- Row 9: "Build reports" ----------------------------------------------
$ws.Range("D9").NumberFormat = "d-mmm-yy"
$ws.Range("D9").Value = 44308
$ws.Range("E9").Value = "4hr"
$ws.Range("E9").HorizontalAlignment = $xlCenter
$ws.Range("E9").VerticalAlignment = $xlCenter
$ws.Range("F9").Value = "completed"
$ws.Range("F9").HorizontalAlignment = $xlCenter

# --- Remarks (G column), entered after the Status/Effort columns ---------
$ws.Range("G8").Value = "Learned powerBi tool and pull the excel sheet dataset to the powerBI"
$ws.Range("G8").VerticalAlignment = $xlTop
$ws.Range("G8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30

$ws.Range("G9").Value = "manually build the E-R diagram and relatioship between the tables"
$ws.Range("G9").VerticalAlignment = $xlTop
$ws.Range("G9").WrapText = $true

# --- Row 10: "Build advanced reports" ------------------------------------
$ws.Range("D10").NumberFormat = "d-mmm-yy"
$ws.Range("D10").Value = 44313
$ws.Range("E10").Value = "2hr"
$ws.Range("E10").HorizontalAlignment = $xlCenter
$ws.Range("E10").VerticalAlignment = $xlCenter
$ws.Range("F10").Value = "completed"
$ws.Range("F10").HorizontalAlignment = $xlCenter
$ws.Range("G10").Value = "understood different graph and explored tooltip, drill down"
$ws.Range("G10").VerticalAlignment = $xlTop
$ws.Range("G10").WrapText = $true

# --- Cosmetic tweaks seen in the diff -------------------------------------
# Column G was widened slightly to better fit the new remarks text.
$ws.Columns.Item(7).ColumnWidth = 65.25

# Cursor/selection ended on G10 after entering the last remark.
$ws.Range("G10").Select() | Out-Null
